# Applies the weekly fruit/vegetable price update (rotates rows 4-13 data
# for Repollo / Agrícola del Norte S.A. de Arica), per the commit
# "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns D (Fecha), I (Calidad), J (Volumen),
# K (Precio mínimo), L (Precio máximo), M (Precio promedio ponderado),
# P (Precio $/Kg) for rows 4 through 13.
$rows = @{
    4  = @{ D = 44229; I = "Segunda"; J = 760;  K = 550; L = 600; M = 575; P = 575 }
    5  = @{ D = 44245; I = "Primera"; J = 800;  K = 850; L = 900; M = 875; P = 875 }
    6  = @{ D = 44245; I = "Segunda"; J = 1000; K = 750; L = 800; M = 775; P = 775 }
    7  = @{ D = 44201; I = "Segunda"; J = 500;  K = 800; L = 900; M = 850; P = 850 }
    8  = @{ D = 44174; I = "Segunda"; J = 800;  K = 450; L = 500; M = 475; P = 475 }
    9  = @{ D = 44174; I = "Tercera"; J = 1200; K = 250; L = 350; M = 300; P = 300 }
    10 = @{ D = 44267; I = "Tercera"; J = 400;  K = 500; L = 600; M = 550; P = 550 }
    11 = @{ D = 44210; I = "Segunda"; J = 900;  K = 600; L = 700; M = 650; P = 650 }
    12 = @{ D = 44278; I = "Segunda"; J = 700;  K = 600; L = 700; M = 650; P = 650 }
    13 = @{ D = 44278; I = "Tercera"; J = 400;  K = 500; L = 600; M = 550; P = 550 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("I$r").Value = $vals.I
    $ws.Range("J$r").Value = $vals.J
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("P$r").Value = $vals.P
}
